# Inserts a new data row at row 84 (pushing the existing rows 84..124 down
# to 85..125) and populates the new row with the record described by the
# commit's diff. This mirrors the target XML diff exactly: dimension grows
# from A1:R124 to A1:R125, and every previously existing row below 83 keeps
# its original data but moves one row further down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 84..124 down to 85..125, creating a blank row 84.
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new record.
$ws.Range("A84").Value = 5
$ws.Range("B84").Value = "Macroferia Regional de Talca"
$ws.Range("C84").Value = "Maule"
$ws.Range("D84").Value = 45141
$ws.Range("E84").Value = 7
$ws.Range("F84").Value = 100112013
$ws.Range("G84").Value = "Alcachofa"
$ws.Range("H84").Value = "Madrigal"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 300
$ws.Range("K84").Value = 14000
$ws.Range("L84").Value = 14000
$ws.Range("M84").Value = 14000
$ws.Range("N84").Value = "$/caja 40 unidades"
$ws.Range("O84").Value = "Provincia del Elquí"
$ws.Range("P84").Value = 350
$ws.Range("Q84").Value = 40
$ws.Range("R84").Value = "Hortaliza"

Write-Host "Inserted new row 84; sheet now spans $($ws.UsedRange.Rows.Count) rows."
